# Absenzenliste-Template: widen the "Name"/"Nachname" column by 2mm
# (1418 -> 1548 dxa) and narrow the "Vorname" column by 2mm
# (1418 -> 1288 dxa). Both columns live in the single table in the
# document (column 3 = Name, column 4 = Vorname); Word propagates the
# new preferred width to every cell of the column (and to <w:tblGrid>).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$nameCol = $t.Columns.Item(3)
$nameCol.Width = 1548 / 20.0

$vornameCol = $t.Columns.Item(4)
$vornameCol.Width = 1288 / 20.0
